# Updates cryptos list figures (Price / Volume(1h) columns) and restores the
# Hedera/Monero row order+values, matching the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.079.56"
$ws.Range("E2").Value = "  +4.27%  "
$ws.Range("D3").Value = "2.346.70"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.Value = "'314.77"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "
$c = $ws.Range("D6")
$c.Value = "'109.23"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.86%  "
$c = $ws.Range("D7")
$c.Value = "'0.630"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("E8").Value = "  +0.00%  "
$c = $ws.Range("D9")
$c.Value = "'0.623"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.36%  "
$c = $ws.Range("D10")
$c.Value = "'41.89"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +8.32%  "
$c = $ws.Range("D11")
$c.Value = "'0.0919"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.25%  "
$c = $ws.Range("D12")
$c.Value = "'8.59"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.88%  "
$ws.Range("E13").Value = "  +4.53%  "
$ws.Range("E14").Value = "  +0.11%  "
$c = $ws.Range("D15")
$c.Value = "'15.50"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").Value = "2.703.24"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "2.346.09"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "44.292.27"
$ws.Range("E18").Value = "  +4.38%  "
$ws.Range("E19").Value = "  +5.11%  "
$ws.Range("E20").Value = "  +3.09%  "
$c = $ws.Range("D21")
$c.Value = "'12.98"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.34%  "
$c = $ws.Range("D22")
$c.Value = "'74.73"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "
$c = $ws.Range("D23")
$c.Value = "'3.50"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$c = $ws.Range("D24")
$c.Value = "'266.45"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.23%  "
$c = $ws.Range("D25")
$c.Value = "'2.27"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.44%  "
$ws.Range("E26").Value = "  -0.63%  "
$c = $ws.Range("D27")
$c.Value = "'7.64"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +10.06%  "
$c = $ws.Range("D28")
$c.Value = "'11.19"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.13%  "
$ws.Range("E29").Value = "  +2.15%  "
$c = $ws.Range("D30")
$c.Value = "'39.66"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +10.47%  "
$c = $ws.Range("D31")
$c.Value = "'22.50"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D32")
$c.Value = "'169.39"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D33")
$c.Value = "'0.0920"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.97%  "
$c = $ws.Range("D34")
$c.Value = "'2.87"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +10.24%  "
$c = $ws.Range("D35")
$c.Value = "'0.132"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("E37").Value = "  +5.88%  "
$c = $ws.Range("D38")
$c.Value = "'0.0366"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.51%  "
$c = $ws.Range("D40")
$c.Value = "'3.80"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  +10.21%  "
$c = $ws.Range("D42")
$c.Value = "'104.03"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.34%  "
$c = $ws.Range("D43")
$c.Value = "'13.99"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +17.71%  "
$c = $ws.Range("D44")
$c.Value = "'0.239"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.75%  "
$c = $ws.Range("D45")
$c.Value = "'71.28"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("E46").Value = "  +0.08%  "
$c = $ws.Range("D47")
$c.Value = "'115.43"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.10%  "
$c = $ws.Range("D48")
$c.Value = "'77.98"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").Value = "1.656.67"
$ws.Range("E49").Value = "  -2.46%  "
$c = $ws.Range("D50")
$c.Value = "'9.04"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.84%  "
$c = $ws.Range("D51")
$c.Value = "'0.217"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +17.04%  "
